$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '27.514.73'
$ws.Cells.Item(2, 5).Value = '  +2.49%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.870.55'
$ws.Cells.Item(3, 5).Value = '  +3.02%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.005'
$ws.Cells.Item(4, 5).Value = '  -0.35%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '316.65'
$ws.Cells.Item(5, 5).Value = '  +2.71%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.004'
$ws.Cells.Item(6, 5).Value = '  -0.33%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4670'
$ws.Cells.Item(7, 5).Value = '  +1.31%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3748'
$ws.Cells.Item(8, 5).Value = '  +3.05%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07410'
$ws.Cells.Item(9, 5).Value = '  +2.74%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.8915'
$ws.Cells.Item(10, 5).Value = '  +4.02%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07966'
$ws.Cells.Item(11, 5).Value = '  +6.21%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '20.11'
$ws.Cells.Item(12, 5).Value = '  +2.12%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.815.83'
$ws.Cells.Item(13, 5).Value = '  +1.27%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.448'
$ws.Cells.Item(14, 5).Value = '  +2.52%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.615'
$ws.Cells.Item(15, 5).Value = '  +1.88%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '92.98'
$ws.Cells.Item(16, 5).Value = '  +1.71%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -0.33%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +4.86%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.004'
$ws.Cells.Item(19, 5).Value = '  -0.28%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.98'
$ws.Cells.Item(20, 5).Value = '  +4.12%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '27.530.37'
$ws.Cells.Item(21, 5).Value = '  +3.05%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.202'
$ws.Cells.Item(22, 5).Value = '  +1.48%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.61'
$ws.Cells.Item(23, 5).Value = '  +1.20%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.078.40'
$ws.Cells.Item(24, 5).Value = '  +4.77%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '153.13'
$ws.Cells.Item(25, 5).Value = '  +1.17%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.885'
$ws.Cells.Item(26, 5).Value = '  +2.16%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '18.63'

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.101'
$ws.Cells.Item(28, 5).Value = '  +1.46%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.181'
$ws.Cells.Item(29, 5).Value = '  +1.95%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '117.63'
$ws.Cells.Item(30, 5).Value = '  +2.13%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08939'
$ws.Cells.Item(31, 5).Value = '  +1.05%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.7566'
$ws.Cells.Item(32, 5).Value = '  +6.41%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.015'
$ws.Cells.Item(33, 5).Value = '  +2.25%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.164'
$ws.Cells.Item(34, 5).Value = '  +3.28%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.509'
$ws.Cells.Item(35, 5).Value = '  +2.39%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.674'
$ws.Cells.Item(36, 5).Value = '  +11.24%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.01973'
$ws.Cells.Item(37, 5).Value = '  +3.16%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Hedera'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.05320'
$ws.Cells.Item(38, 5).Value = '  +1.70%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.083'
$ws.Cells.Item(39, 5).Value = '  +1.07%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.992'
$ws.Cells.Item(40, 5).Value = '  +2.64%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '7.232'
$ws.Cells.Item(41, 5).Value = '  +1.63%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.5260'
$ws.Cells.Item(42, 5).Value = '  +2.61%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1653'
$ws.Cells.Item(43, 5).Value = '  +2.14%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '8.385'

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.4921'
$ws.Cells.Item(45, 5).Value = '  +3.03%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '10.41'
$ws.Cells.Item(46, 5).Value = '  +3.42%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.681'
$ws.Cells.Item(47, 5).Value = '  +4.21%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'PaxDollar'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.005'
$ws.Cells.Item(48, 5).Value = '  -0.29%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Quant'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '104.13'
$ws.Cells.Item(49, 5).Value = '  +1.30%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.06266'
$ws.Cells.Item(50, 5).Value = '  +0.21%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '66.01'
$ws.Cells.Item(51, 5).Value = '  +3.41%  '
